# The "Förändrad" (Changed) column C was bumped by one day (45178 -> 45179)
# for every data row (rows 2 through 210) on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C210").Value = 45179
